# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values on the zh-cn and de-de
# report sheets to reflect the freshly (re)generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-14 01:22:07"
$wsZhCn.Range("H2").Value = "2016-03-14 01:22:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-14 01:22:11"
$wsDeDe.Range("H2").Value = "2016-03-14 01:22:31"
